# "Turned all parsim into functions"
# Update the computed mean mass flow rate values (column B) and the
# resulting auto-fit column widths for columns A:C on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New mean mass flow rate values (re-computed from functions instead of
# the previously parsed/static numbers).
$ws.Range("B2").Value = 0.056162696147666641
$ws.Range("B3").Value = 0.48579507895945706
$ws.Range("B4").Value = 0.07262684780040872
$ws.Range("B5").Value = 0.041585799737799654
$ws.Range("B6").Value = 0.2206948521042785
$ws.Range("B7").Value = 0.077034421117002455
$ws.Range("B8").Value = 0.038972319652292191
$ws.Range("B9").Value = 0.33260624224759788

# Column widths shrank slightly to better fit the new values/labels.
$ws.Columns.Item(1).ColumnWidth = 23.166666666666668
$ws.Columns.Item(2).ColumnWidth = 17
$ws.Columns.Item(3).ColumnWidth = 3.8333333333333335
